# Adds a new data-collection wave (28. 12. 2021) to both sheets of the
# workbook, appending a new trailing column (AL on "data", AK on "pocetR"),
# and refreshes the "aktualizace" (updated-on) date in each sheet's footer
# row from 8. 12. 2021 to 6. 1. 2022.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "data" -- percentages, new column AL, header date 28. 12. 2021
# -----------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# New header cell, formatted like its left neighbour (AK1).
$wsData.Range("AK1").Copy()
$wsData.Range("AL1").PasteSpecial(-4122)   # xlPasteFormats
$wsData.Range("AL1").Value = "28. 12. 2021"

# New data values for rows 2..45 (column AL), in row order.
$alValues = @(
    0.2, 0.13, 0.43, 0.3, 0.14, 0.21, 0.18, 0.19, 0.2, 0.2,
    0.15, 0.32, 0.19, 0.2, 0.22, 0.17, 0.26, 0.25, 0.18, 0.11,
    0.15, 0.16, 0.46, 0.42, 0.14, 0.08, 0.14, 0.12, 0.12, 0.11,
    0.13, 0.12, 0.19, 0.11, 0.16, 0.12, 0.06, 0.26, 0.14, 0.11,
    0.05, 0.07, 0.16, 0.27
)

for ($i = 0; $i -lt $alValues.Length; $i++) {
    $row = $i + 2
    $wsData.Cells.Item($row, 38).Value = $alValues[$i]
}

# Footer row (46): bump the "aktualizace" date.
$wsData.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 6. 1. 2022"

# -----------------------------------------------------------------
# Sheet "pocetR" -- sample sizes, new column AK, header date 28. 12. 2021
# -----------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

# New header cell, formatted like its left neighbour (AJ1).
$wsPocet.Range("AJ1").Copy()
$wsPocet.Range("AK1").PasteSpecial(-4122)   # xlPasteFormats
$wsPocet.Range("AK1").Value = "28. 12. 2021"

# New data values for rows 2..23 (column AK), in row order.
$akValues = @(
    1767, 192, 333, 1242, 851, 157, 505, 254, 817, 143,
    112, 695, 818, 606, 343, 190, 640, 596, 244, 522,
    290, 182
)

for ($i = 0; $i -lt $akValues.Length; $i++) {
    $row = $i + 2
    $wsPocet.Cells.Item($row, 37).Value = $akValues[$i]
}

# Trailing blank cell in the footer row, matching the rest of row 24.
$wsPocet.Range("AK24").Value = ""

# Footer row (24): bump the "aktualizace" date.
$wsPocet.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 1. 2022"
